$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2026-01-02 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-01-03 Saturday", 2)

# Update the multiplication problems in the table, cell by cell, using
# (row, column) coordinates so that duplicate expressions (e.g. "670x2=")
# are each replaced with their correct, distinct target value.
$table = $d.Tables.Item(1)

$cellUpdates = @(
    @{ Row = 1;  Col = 1; New = "148×4=" },
    @{ Row = 1;  Col = 2; New = "361×5=" },
    @{ Row = 1;  Col = 3; New = "738×8=" },
    @{ Row = 1;  Col = 4; New = "183×6=" },
    @{ Row = 1;  Col = 5; New = "970×3=" },

    @{ Row = 5;  Col = 1; New = "497×5=" },
    @{ Row = 5;  Col = 2; New = "297×9=" },
    @{ Row = 5;  Col = 3; New = "828×3=" },
    @{ Row = 5;  Col = 4; New = "147×7=" },
    @{ Row = 5;  Col = 5; New = "279×9=" },

    @{ Row = 10; Col = 1; New = "387×9=" },
    @{ Row = 10; Col = 2; New = "244×8=" },
    @{ Row = 10; Col = 3; New = "466×8=" },
    @{ Row = 10; Col = 4; New = "640×6=" },
    @{ Row = 10; Col = 5; New = "392×5=" },

    @{ Row = 15; Col = 1; New = "755×3=" },
    @{ Row = 15; Col = 2; New = "880×2=" },
    @{ Row = 15; Col = 3; New = "182×7=" },
    @{ Row = 15; Col = 4; New = "318×6=" },
    @{ Row = 15; Col = 5; New = "294×6=" },

    @{ Row = 20; Col = 1; New = "379×3=" },
    @{ Row = 20; Col = 2; New = "770×6=" },
    @{ Row = 20; Col = 3; New = "631×2=" },
    @{ Row = 20; Col = 4; New = "670×8=" },
    @{ Row = 20; Col = 5; New = "561×7=" }
)

foreach ($update in $cellUpdates) {
    $cell = $table.Cell($update.Row, $update.Col)
    $cellRange = $cell.Range
    $cellRange.MoveEnd(12, -1) | Out-Null
    $cellRange.Text = $update.New
}
